# "added true and false predicates to the alloy specifications that we
# generate" -- the spec document has two small operation tables ("init"
# and "save"). Each table has a row whose second cell held the
# placeholder "some PID" and, for the "save" table, two more rows whose
# second cell held the placeholder formula "{ p : PID | p != p }".
# Those placeholders become the literal Alloy predicates "true" / "none".
#
# Table objects are re-fetched right before each edit (rather than
# reused across edits) because this host's Table/Cell handles can go
# stale once an earlier edit has changed the document.

$d = $word.ActiveDocument

function Set-CellText {
    param($TableIndex, $RowIndex, $ColIndex, $NewText)

    $table = $d.Tables.Item($TableIndex)
    $cell = $table.Cell($RowIndex, $ColIndex)
    $range = $cell.Range
    # Trim the trailing end-of-cell marker so we only replace the text.
    $range.End = $range.End - 1
    $range.Text = $NewText
}

# Table 4 -- the "init" operation table: blank/"some PID" row -> "true".
Set-CellText 4 3 2 "true"

# Table 7 -- the "save" operation table.
# blank/"some PID" row -> "true".
Set-CellText 7 3 2 "true"
# "toAdd" row's formula -> "none".
Set-CellText 7 4 2 "none"
# "toDelete" row's formula -> "none".
Set-CellText 7 5 2 "none"
